$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "cdf" column actually plotted an empirical CDF, so rename the header.
$ws.Range("H1").Value = "ecdf"

# Recompute the empirical CDF as rank / sample_size (previously it used
# rank / (sample_size + 1), a plotting-position convention). There are 95
# observations (rows 2-96), with the ascending rank already available in
# column E.
$sampleSize = 95
for ($row = 2; $row -le 96; $row++) {
    $rank = $ws.Cells.Item($row, 5).Value2
    $ws.Cells.Item($row, 8).Value = $rank / $sampleSize
}
